$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1213.3334
$ws.Range("I19").Value = 1213.3334
$ws.Range("K19").Value = 1213.3334
$ws.Range("M19").Value = -1038.3334

$ws.Range("H38").Value = 1250127.1
$ws.Range("I38").Value = 1666749.6
$ws.Range("K38").Value = 5000248.800000001
$ws.Range("M38").Value = -4999876.800000001

$ws.Range("H64").Value = 5979.4
$ws.Range("I64").Value = 6224.5
$ws.Range("K64").Value = 6224.5
$ws.Range("M64").Value = -5976.5

$ws.Range("H67").Value = 5979.4
$ws.Range("I67").Value = 6224.5
$ws.Range("K67").Value = 6224.5
$ws.Range("M67").Value = -5366.5

$ws.Range("H80").Value = 4780.5
$ws.Range("I80").Value = 4751
$ws.Range("J80").Value = 4790.3335
$ws.Range("K80").Value = 14253
$ws.Range("L80").Value = 14371.0005
$ws.Range("M80").Value = -13255
$ws.Range("N80").Value = -16367.0005

$ws.Range("H83").Value = 4780.5
$ws.Range("I83").Value = 4751
$ws.Range("J83").Value = 4790.3335
$ws.Range("K83").Value = 42759
$ws.Range("L83").Value = 43113.0015
$ws.Range("M83").Value = -37767
$ws.Range("N83").Value = -53097.0015

$ws.Range("H97").Value = 1450
$ws.Range("J97").Value = 2200
$ws.Range("L97").Value = 6600
$ws.Range("N97").Value = -7592

$ws.Range("H98").Value = 683.7778
$ws.Range("I98").Value = 683.7778
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 683.7778
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 814.2222
$ws.Range("N98").ClearContents()

$ws.Range("H100").Value = 2474.5
$ws.Range("I100").Value = 2312
$ws.Range("K100").Value = 2312
$ws.Range("M100").Value = -1771

$ws.Range("H112").Value = 3982.2856
$ws.Range("J112").Value = 3982.2856
$ws.Range("L112").Value = 11946.8568
$ws.Range("N112").Value = -14162.8568

$ws.Range("H116").Value = 4674.75
$ws.Range("I116").Value = 4674.75
$ws.Range("K116").Value = 4674.75
$ws.Range("M116").Value = -1232.75

$ws.Range("H122").Value = 683.7778
$ws.Range("I122").Value = 683.7778
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2051.3334
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 398.6666
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4202
$ws.Range("I32").Value = 3142.2
$ws.Range("K32").Value = 3142.2
$ws.Range("M32").Value = -2855.2

$ws.Range("H45").Value = 2249.6667
$ws.Range("I45").Value = 2318.375
$ws.Range("J45").Value = 1700
$ws.Range("K45").Value = 2318.375
$ws.Range("L45").Value = 1700
$ws.Range("M45").Value = -1941.375
$ws.Range("N45").Value = -2454

$ws.Range("H132").Value = 7653.5
$ws.Range("I132").Value = 7392.778
$ws.Range("K132").Value = 22178.334
$ws.Range("M132").Value = -19648.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2149
$ws.Range("I62").Value = 2149
$ws.Range("K62").Value = 2149
$ws.Range("M62").Value = -1525

$ws.Range("H65").Value = 2149
$ws.Range("I65").Value = 2149
$ws.Range("K65").Value = 10745
$ws.Range("M65").Value = -7625

$ws.Range("H86").Value = 7884.4443
$ws.Range("I86").Value = 7490
$ws.Range("K86").Value = 7490
$ws.Range("M86").Value = -6367

$ws.Range("H89").Value = 7884.4443
$ws.Range("I89").Value = 7490
$ws.Range("K89").Value = 37450
$ws.Range("M89").Value = -31834

$ws.Range("H109").Value = 99999.89999999999
$ws.Range("J109").Value = 99999.89999999999
$ws.Range("L109").Value = 99999.89999999999
$ws.Range("N109").Value = -102079.9

$ws.Range("H122").Value = 6089.778
$ws.Range("I122").Value = 7632.5
$ws.Range("K122").Value = 22897.5
$ws.Range("M122").Value = -20447.5

$ws.Range("H134").Value = 3583.3333
$ws.Range("I134").Value = 3468.75
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 10406.25
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -7871.25
$ws.Range("N134").Value = -18570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 94.28570999999999
$ws.Range("J2").Value = 107.125
$ws.Range("L2").Value = 642.75
$ws.Range("N2").Value = -868.75

$ws.Range("H38").Value = 295
$ws.Range("I38").Value = 363.33334
$ws.Range("K38").Value = 1090.00002
$ws.Range("M38").Value = -743.0000199999999

$ws.Range("H131").Value = 61957.69
$ws.Range("I131").Value = 1628
$ws.Range("J131").Value = 113668.86
$ws.Range("K131").Value = 4884
$ws.Range("L131").Value = 341006.58
$ws.Range("M131").Value = 156
$ws.Range("N131").Value = -351086.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 35000
$ws.Range("J69").Value = 35000
$ws.Range("L69").Value = 35000
$ws.Range("N69").Value = -36498

$ws.Range("H70").Value = 9568.154
$ws.Range("I70").Value = 9873.083000000001
$ws.Range("J70").Value = 5909
$ws.Range("K70").Value = 9873.083000000001
$ws.Range("L70").Value = 5909
$ws.Range("M70").Value = -9603.083000000001
$ws.Range("N70").Value = -6449

$ws.Range("H72").Value = 35000
$ws.Range("J72").Value = 35000
$ws.Range("L72").Value = 105000
$ws.Range("N72").Value = -112488

$ws.Range("H73").Value = 9568.154
$ws.Range("I73").Value = 9873.083000000001
$ws.Range("J73").Value = 5909
$ws.Range("K73").Value = 9873.083000000001
$ws.Range("L73").Value = 5909
$ws.Range("M73").Value = -8937.083000000001
$ws.Range("N73").Value = -7781

$ws.Range("H122").Value = 2630.6316
$ws.Range("I122").Value = 2146.2
$ws.Range("K122").Value = 6438.599999999999
$ws.Range("M122").Value = -3988.599999999999

$ws.Range("H132").Value = 2899.4
$ws.Range("I132").Value = 2910.6667
$ws.Range("J132").Value = 2798
$ws.Range("K132").Value = 8732.000100000001
$ws.Range("L132").Value = 8394
$ws.Range("M132").Value = -6202.000100000001
$ws.Range("N132").Value = -13454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 3666.6667
$ws.Range("I56").Value = 2500
$ws.Range("J56").Value = 4250
$ws.Range("K56").Value = 2500
$ws.Range("L56").Value = 4250
$ws.Range("M56").Value = -1809
$ws.Range("N56").Value = -5632

$ws.Range("H61").Value = 1913.3334
$ws.Range("I61").Value = 1902.5
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1902.5
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1700.5
$ws.Range("N61").Value = -2404

$ws.Range("H94").Value = 59500
$ws.Range("J94").Value = 59500
$ws.Range("L94").Value = 59500
$ws.Range("N94").Value = -60852

$ws.Range("H113").Value = 1913.3334
$ws.Range("I113").Value = 1902.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1902.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 267.5
$ws.Range("N113").Value = -6340

$ws.Range("H122").Value = 3829
$ws.Range("I122").Value = 3358
$ws.Range("K122").Value = 10074
$ws.Range("M122").Value = -7624

$ws.Range("H132").Value = 6066.6665
$ws.Range("I132").Value = 3200
$ws.Range("K132").Value = 9600
$ws.Range("M132").Value = -7070

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 35891
$ws.Range("I58").Value = 6836.5
$ws.Range("K58").Value = 6836.5
$ws.Range("M58").Value = -6528.5

$ws.Range("H130").Value = 56497.5
$ws.Range("J130").Value = 56497.5
$ws.Range("L130").Value = 56497.5
$ws.Range("N130").Value = -66537.5

$ws.Range("H132").Value = 3381.8572
$ws.Range("I132").Value = 3001.1177
$ws.Range("K132").Value = 9003.3531
$ws.Range("M132").Value = -6473.3531
